# aging_stock_summary_copy.xlsx data refresh
# Updates a handful of per-SKU aging-stock quantities on Sheet1 (rows
# 15, 19, 21, 22, 24) to reflect the latest inventory pull, and clears
# out cells whose counts are no longer present in the refreshed extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 (91 - 180 Days / Losectil)
$ws.Range("AE15").Value = 33

# Row 19 (91 - 180 Days / Losectil) - CTGSKF count removed in refresh
$ws.Range("I19").Value = $null

# Row 21 (91 - 180 Days / Panoral)
$ws.Range("G21").Value = 16
$ws.Range("AE21").Value = 45

# Row 22 (181 - 210 Days / Remivir)
$ws.Range("I22").Value = 30
$ws.Range("P22").Value = 14
$ws.Range("V22").Value = $null
$ws.Range("Z22").Value = 105
$ws.Range("AD22").Value = $null
$ws.Range("AE22").Value = 8

# Row 24 (181 - 210 Days / Flucoder)
$ws.Range("V24").Value = 20
